$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: overwrite C22 (was "P34"), add D22 / E22 ---
$ws.Range("C22").Value = "0.125 × 5 ="

$ws.Range("D22").Value = "1.5 ÷ 40 ="
$ws.Range("D22").VerticalAlignment = -4108

$ws.Range("E22").Value = "4.8 ÷ 60 ="
$ws.Range("E22").VerticalAlignment = -4108

# --- Row 23 (new) ---
$ws.Range("A23").Value = "8.5 ÷ 1.7 ="
$ws.Range("B23").Value = "8 ÷ 0.5 ="
$ws.Range("C23").Value = "0.92 ÷ 0.4 ="
$ws.Range("D23").Value = "0.25 × 60 ="
$ws.Range("E23").Value = "0.8 ÷ 0.01 ="
$ws.Range("A23:E23").VerticalAlignment = -4108
$ws.Rows.Item(23).RowHeight = 23

# --- Row 24 (new) ---
$ws.Range("A24").Value = "3.2 ÷ 0.2 ="
$ws.Range("B24").Value = "7.2 ÷ 0.4 ="
$ws.Range("C24").Value = "3.7 ÷ 0.5 ="
$ws.Range("D24").Value = "2.6 × 4 ="
$ws.Range("E24").Value = "1.3 × 0.5 ="
$ws.Range("A24:E24").VerticalAlignment = -4108
$ws.Rows.Item(24).RowHeight = 23

# --- Row 25 (new) ---
$ws.Range("A25").Value = "1.2 × 0.4 ="
$ws.Range("B25").Value = "0.52b + 0.8b ="
$ws.Range("C25").Value = "7.6 ÷ 0.2 ="
$ws.Range("D25").Value = "0.75 ÷ 0.25 ="
$ws.Range("E25").Value = "0.92 ÷ 0.4 ="
$ws.Range("A25:E25").VerticalAlignment = -4108
$ws.Rows.Item(25).RowHeight = 23

# --- Row 26 (new) ---
$ws.Range("A26").Value = "8.6χ -7χ = 32"
$ws.Range("B26").Value = " χ ="
$ws.Range("C26").Value = "9χ -4χ = 7.8"
$ws.Range("D26").Value = " χ ="
$ws.Range("E26").Value = "3.6 × 0.5 ="
$ws.Range("A26:E26").VerticalAlignment = -4108
$ws.Rows.Item(26).RowHeight = 23
# split "8.6" / "χ -7χ = 32" into two runs by touching the font of the tail
$ws.Range("A26").Characters(4, 10).Font.Name = "Calibri"
# give C26 its own explicit (non-theme) font, matching the source file
$ws.Range("C26").Font.Name = "Calibri"

# --- Row 27 (new) ---
$ws.Range("A27").Value = "0.4χ  = 16.4"
$ws.Range("B27").Value = " χ ="
$ws.Range("C27").Value = "3.4 ÷ 68 ="
$ws.Range("D27").Value = "P49"
$ws.Range("E27").Value = ""
$ws.Range("A27:E27").VerticalAlignment = -4108
$ws.Rows.Item(27).RowHeight = 24
# split "0.4" / "χ  = 16.4" into two runs
$ws.Range("A27").Characters(4, 9).Font.Name = "Calibri"

# --- view / selection updates ---
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D28").Select()
